$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Insert a new column at H: the old "Notes" column (H) shifts right to I,
# and the new column H becomes "Arduino Pin".
$ws.Columns("H").Insert()

# Header row
$ws.Range("H2").Value = "Arduino Pin"

# New Arduino Pin numbers for two rows (Beaglebone Status Lights relay rows)
$ws.Range("H14").Value = 39
$ws.Range("H17").Value = 38

# Column C previously held numeric pin-position placeholders (1/2/3); they
# are replaced with the actual wire colors used for those pins.
$ws.Range("C9").Value = "blue"
$ws.Range("C10").Value = "black"
$ws.Range("C11").Value = "brown"

$ws.Range("C13").Value = "blue"
$ws.Range("C14").Value = "brown"
$ws.Range("C15").Value = "black"

$ws.Range("C16").Value = "blue"
$ws.Range("C17").Value = "brown"
$ws.Range("C18").Value = "black"

$ws.Range("C19").Value = "blue"
$ws.Range("C20").Value = "brown"
$ws.Range("C21").Value = "black"

$ws.Range("C22").Value = "blue"
$ws.Range("C23").Value = "brown"
$ws.Range("C24").Value = "black"

$ws.Range("C45").Value = "blue"
$ws.Range("C46").Value = "black"

# Match the saved view/selection state from the edit.
$ws.Range("I15").Select()
$excel.ActiveWindow.ScrollRow = 4
